$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the contents of a whole paragraph (found by matching its
# current text) with a hand-built run of OOXML (used to splice in
# <w:proofErr> spell/grammar-check markers that Find/Replace can't produce).
# ---------------------------------------------------------------------------
function Set-ParagraphXml($matchText, $innerXml) {
    $paras = $d.Paragraphs
    $count = $paras.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        $t = $p.Range.Text.TrimEnd("`r")
        if ($t -eq $matchText) {
            $r = $p.Range
            $xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $r.InsertXML($xml)
            return $true
        }
    }
    return $false
}

# ---------------------------------------------------------------------------
# 1) "*it's important ... Optrex C-51505 ..." paragraph:
#    - wrap "it's" in gramStart/gramEnd proofErr markers
#    - split "a Optrex C-51505" so "Optrex" is wrapped in spellStart/spellEnd
# ---------------------------------------------------------------------------
$para1Old = "*it’s important to note that we are using a Optrex C-51505 as other common LCD shields can have "
$para1Xml = '<w:p w14:paraId="24C5DED1" w14:textId="1C5CF1DB" w:rsidR="001C5633" w:rsidRDefault="00881216" w:rsidP="001C5633"><w:r><w:t>*</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="001C5633"><w:t>it’s</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> important to note that we are using a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Optrex</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> C-51505 as other common LCD shields ca</w:t></w:r><w:r w:rsidR="001C5633"><w:t xml:space="preserve">n have </w:t></w:r></w:p>'
Set-ParagraphXml $para1Old $para1Xml

# ---------------------------------------------------------------------------
# 2) "We used the SimpleDHT library ..." paragraph:
#    - split off "SimpleDHT" and wrap it in spellStart/spellEnd proofErr markers
# ---------------------------------------------------------------------------
$para2Old = "We used the SimpleDHT library to read input from the DHT11 sensor, which measures temperature and humidity. With these values, we were able to calculate the heat index, which determines the level of comfort people face depending on the temperature and humidity of their environment. We used the MQ2 library to find the average smoke level from multiple inputs from the MQ2 sensor. Finally, we used the MQ9 library to find the average carbon monoxide level from multiple inputs from the MQ9 sensor. After determining heat index, smoke level, and carbon monoxide level, we determined the different thresholds for each value. We then set the frequency of the LED to represent the severity of each value. The higher the frequency of the LED, the more severe the value was determined to be. "
$para2Xml = '<w:p w14:paraId="6FDF0821" w14:textId="75484F85" w:rsidR="00375900" w:rsidRDefault="00FF2453" w:rsidP="00183133"><w:pPr><w:ind w:firstLine="0"/></w:pPr><w:r><w:t xml:space="preserve">We used the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SimpleDHT</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> library to read input from the DHT11 sensor, which measures temperature and humidity. With these values, we were able to calculate the heat index, which determines the level of comfort people face depending on the temperature and humidity of their environment. We used the MQ2 library to find the average smoke level from multiple inputs from the MQ2 sensor. Finally, we used the MQ9 library to find the average carbon monoxide level from multiple inputs from the MQ9 sensor. After determining heat index, smoke level, and carbon monoxide level, we determined the different thresholds for each v</w:t></w:r><w:r w:rsidR="002B2F17"><w:t xml:space="preserve">alue. </w:t></w:r><w:r><w:t xml:space="preserve">We then set the frequency of the LED to represent the severity of each value. The higher the frequency of the LED, the more severe the value was determined to be. </w:t></w:r></w:p>'
Set-ParagraphXml $para2Old $para2Xml

# ---------------------------------------------------------------------------
# 3) Remove the first of the run of empty underline-formatted paragraphs
#    that follows the "Output: ... based on its severity." paragraph.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if ($t -like "*based on its severity.*") {
        $empty = $paras.Item($i + 1)
        $empty.Range.Delete()
        break
    }
}
